{"js": "// Locate the three paragraphs that form the \"python model_main...\" /\n// blank / \"python object_detection/model_main...\" block near the end of\n// the document, then replace that whole range in a single shot with the\n// post-edit OOXML: the first two runs get double-strikethrough (the\n// commands are superseded), a new \"legacy/train.py\" command is added\n// together with an introductory sentence, and the now-redundant trailing\n// \"_GoBack\" bookmark moves onto the new final paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet startIndex = -1;\nlet endIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"python model_main.py\") !== -1 && t.indexOf(\"mob_train_dir\") !== -1) {\n    startIndex = i;\n  }\n  if (t.indexOf(\"python object_detection/model_main.py\") !== -1) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1 || endIndex <= startIndex) {\n  throw new Error(\"Could not locate the target paragraph block\");\n}\n\nconst startRange = items[startIndex].getRange(\"Start\");\nconst endRange = items[endIndex].getRange(\"End\");\nconst targetRange = startRange.expandTo(endRange);\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t>python model_main.py --</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t>model_dir</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> voc/mob_train_dir/ --pipeline_config_path voc/mob.config</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:lastRenderedPageBreak/>\n    <w:t>python object_detection/model_main.py --model_dir object_detection/voc/train_dir/ --pipeline_config_path object_detection/voc/voc.config</w:t>\n  </w:r>\n</w:p>\n<w:p/>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n    <w:t>\u6709\u4e00\u4e2a</w:t>\n  </w:r>\n  <w:r>\n    <w:t>legecy</w:t>\n  </w:r>\n  <w:r>\n    <w:t>\u7684</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:t>train</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:t>python legacy/train.py --train_dir vvv/trainout/ --pipeline_config_path vvv/vv.config</w:t>\n  </w:r>\n  <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n  <w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Locate the three paragraphs that form the \"python model_main...\" /\n# blank / \"python object_detection/model_main...\" block near the end of\n# the document, then replace that whole range in a single shot with the\n# post-edit OOXML: the first two runs get double-strikethrough (the\n# commands are superseded), a new \"legacy/train.py\" command is added\n# together with an introductory sentence, and the now-redundant trailing\n# \"_GoBack\" bookmark moves onto the new final paragraph.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$startIndex = -1\n$endIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*python model_main.py*\" -and $t -like \"*mob_train_dir*\") {\n        $startIndex = $i\n    }\n    if ($t -like \"*python object_detection/model_main.py*\") {\n        $endIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1 -or $endIndex -eq -1) {\n    throw \"Could not locate the target paragraph block\"\n}\n\n$startPara = $d.Paragraphs.Item($startIndex)\n$endPara = $d.Paragraphs.Item($endIndex)\n$targetRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t>python model_main.py --</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t>model_dir</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> voc/mob_train_dir/ --pipeline_config_path voc/mob.config</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n      <w:dstrike/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:dstrike/>\n    </w:rPr>\n    <w:lastRenderedPageBreak/>\n    <w:t>python object_detection/model_main.py --model_dir object_detection/voc/train_dir/ --pipeline_config_path object_detection/voc/voc.config</w:t>\n  </w:r>\n</w:p>\n<w:p/>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n    <w:t>\u6709\u4e00\u4e2a</w:t>\n  </w:r>\n  <w:r>\n    <w:t>legecy</w:t>\n  </w:r>\n  <w:r>\n    <w:t>\u7684</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:t>train</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:rFonts w:hint=\"eastAsia\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:t>python legacy/train.py --train_dir vvv/trainout/ --pipeline_config_path vvv/vv.config</w:t>\n  </w:r>\n  <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n  <w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n\n$targetRange.InsertXML($ooxml) | Out-Null\n"}
